# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "29.174.97"
$ws.Range("E2").Value = "  +0.41%  "
Set-TextCell "D3" "1.829.18"
$ws.Range("E3").Value = "  -0.24%  "
Set-TextCell "D4" "0.9991"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextCell "D5" "242.83"
$ws.Range("E5").Value = "  +0.00%  "
Set-TextCell "D6" "0.6149"
$ws.Range("E6").Value = "  -0.76%  "
Set-TextCell "D7" "1.000"
$ws.Range("E7").Value = "  -0.01%  "
Set-TextCell "D8" "0.07337"
$ws.Range("E8").Value = "  -1.79%  "
Set-TextCell "D9" "0.2907"
Set-TextCell "D10" "23.16"
$ws.Range("E10").Value = "  +0.22%  "
Set-TextCell "D11" "0.07629"
$ws.Range("E11").Value = "  -0.60%  "
Set-TextCell "D12" "1.828.05"
$ws.Range("E12").Value = "  +0.05%  "
Set-TextCell "D13" "4.974"
$ws.Range("E13").Value = "  -0.62%  "
Set-TextCell "D14" "0.6707"
$ws.Range("E14").Value = "  -0.38%  "
Set-TextCell "D15" "82.41"
$ws.Range("E15").Value = "  -0.57%  "
Set-TextCell "D16" "0.000008970"
$ws.Range("E16").Value = "  -2.11%  "
Set-TextCell "D17" "5.843"
$ws.Range("E17").Value = "  -1.21%  "
Set-TextCell "D18" "29.162.74"
$ws.Range("E18").Value = "  +0.33%  "
Set-TextCell "D19" "2.079.38"
$ws.Range("E19").Value = "  -0.08%  "
Set-TextCell "D20" "236.34"
$ws.Range("E20").Value = "  -1.39%  "
Set-TextCell "D21" "12.49"
$ws.Range("E21").Value = "  -1.61%  "
Set-TextCell "D22" "1.000"
$ws.Range("E22").Value = "  -0.03%  "
Set-TextCell "D23" "7.359"
$ws.Range("E23").Value = "  +2.12%  "
Set-TextCell "D24" "1.000"
$ws.Range("E24").Value = "  +0.02%  "
Set-TextCell "D25" "158.63"
$ws.Range("E25").Value = "  -0.42%  "
Set-TextCell "D26" "8.519"
$ws.Range("E26").Value = "  +0.21%  "
Set-TextCell "D27" "0.1385"
$ws.Range("E27").Value = "  -2.99%  "
Set-TextCell "D28" "17.62"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  -0.77%  "
Set-TextCell "D30" "0.05821"
$ws.Range("E30").Value = "  +4.57%  "
Set-TextCell "B31" "Toncoin"
Set-TextCell "C31" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D31" "1.220"
$ws.Range("E31").Value = "  +1.08%  "
Set-TextCell "B32" "InternetComputer(DFINITY)"
Set-TextCell "C32" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D32" "4.074"
$ws.Range("E32").Value = "  -1.13%  "
Set-TextCell "D33" "4.081"
$ws.Range("E33").Value = "  -1.68%  "
Set-TextCell "D34" "1.850"
$ws.Range("E34").Value = "  +0.45%  "
Set-TextCell "D35" "1.134"
$ws.Range("E35").Value = "  -0.59%  "
Set-TextCell "D36" "0.7173"
$ws.Range("E36").Value = "  -2.74%  "
Set-TextCell "D37" "2.615"
$ws.Range("E37").Value = "  -1.59%  "
Set-TextCell "D38" "2.864"
$ws.Range("E38").Value = "  +3.29%  "
Set-TextCell "D39" "1.226.40"
$ws.Range("E39").Value = "  +1.29%  "
Set-TextCell "D40" "0.01759"
$ws.Range("E40").Value = "  -1.27%  "
Set-TextCell "D41" "6.183"
$ws.Range("E41").Value = "  -4.48%  "
Set-TextCell "D42" "0.8993"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +0.11%  "
Set-TextCell "D44" "1.995.11"
$ws.Range("E44").Value = "  +0.86%  "
Set-TextCell "D45" "101.75"
$ws.Range("E45").Value = "  +0.01%  "
Set-TextCell "D46" "65.51"
$ws.Range("E46").Value = "  -0.26%  "
Set-TextCell "D47" "0.5040"
$ws.Range("E47").Value = "  -0.87%  "
Set-TextCell "B48" "EnergySwap"
Set-TextCell "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D48" "9.194"
$ws.Range("E48").Value = "  +0.55%  "
Set-TextCell "B49" "TheSandbox"
Set-TextCell "C49" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D49" "0.4037"
$ws.Range("E49").Value = "  -0.88%  "
Set-TextCell "B50" "BabyDogeCoin"
Set-TextCell "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D50" "0.00000000117"
$ws.Range("E50").Value = "  -1.44%  "
Set-TextCell "D51" "0.1152"
$ws.Range("E51").Value = "  +4.34%  "
